$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.9
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 2.5
$ws.Range("J2").Value = 1.07
$ws.Range("K2").Value = 8.5
$ws.Range("T2").Value = 8.5
$ws.Range("X2").Value = 26
$ws.Range("Z2").Value = 8.5
$ws.Range("AF2").Value = 11
$ws.Range("AH2").Value = 23
# Row 4
$ws.Range("K4").Value = 15
# Row 5
$ws.Range("L5").Value = 1.29
$ws.Range("M5").Value = 3.5
$ws.Range("N5").Value = 1.95
$ws.Range("O5").Value = 1.85
# Row 6
$ws.Range("J6").Value = 1.11
$ws.Range("K6").Value = 6.5
# Row 8
$ws.Range("J8").Value = 1.11
$ws.Range("K8").Value = 6.5
# Row 11
$ws.Range("K11").Value = 7.5
$ws.Range("AA11").Value = 6
$ws.Range("AB11").Value = 17
$ws.Range("AD11").Value = 451
$ws.Range("AE11").Value = 10
$ws.Range("AJ11").Value = 51
# Row 22
$ws.Range("J22").Value = 1.07
$ws.Range("K22").Value = 9
# Row 23
$ws.Range("L23").Value = 1.3
$ws.Range("M23").Value = 3.4
$ws.Range("N23").Value = 2.03
$ws.Range("O23").Value = 1.78
# Row 24
$ws.Range("N24").Value = 2
$ws.Range("O24").Value = 1.8
# Row 25
$ws.Range("N25").Value = 1.57
$ws.Range("O25").Value = 2.35
# Row 28
$ws.Range("G28").Value = 3.4
$ws.Range("I28").Value = 2
$ws.Range("V28").Value = 9.75
$ws.Range("X28").Value = 24
$ws.Range("Y28").Value = 29
$ws.Range("Z28").Value = 9
$ws.Range("AG28").Value = 7.3
# Row 30
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 5
$ws.Range("I30").Value = 1.38
$ws.Range("K30").Value = 21
$ws.Range("R30").Value = 1.67
$ws.Range("S30").Value = 2.1
$ws.Range("W30").Value = 81
$ws.Range("Z30").Value = 21
$ws.Range("AA30").Value = 10
$ws.Range("AB30").Value = 17
$ws.Range("AE30").Value = 10
$ws.Range("AG30").Value = 8.5
$ws.Range("AH30").Value = 10
# Row 31
$ws.Range("G31").Value = 1.85
$ws.Range("I31").Value = 4.1
$ws.Range("P31").Value = 1.33
$ws.Range("Q31").Value = 3.25
$ws.Range("R31").Value = 1.7
$ws.Range("S31").Value = 2.05
$ws.Range("W31").Value = 15
$ws.Range("Y31").Value = 23
$ws.Range("Z31").Value = 12
$ws.Range("AA31").Value = 7
$ws.Range("AB31").Value = 13
$ws.Range("AE31").Value = 13
$ws.Range("AI31").Value = 29
# Row 32
$ws.Range("G32").Value = 13
$ws.Range("I32").Value = 1.21
$ws.Range("L32").Value = 1.22
$ws.Range("M32").Value = 3.45
$ws.Range("N32").Value = 1.65
$ws.Range("O32").Value = 1.98
$ws.Range("R32").Value = 2.4
$ws.Range("S32").Value = 1.44
$ws.Range("T32").Value = 28
$ws.Range("U32").Value = 110
$ws.Range("V32").Value = 45
$ws.Range("W32").Value = 600
$ws.Range("X32").Value = 250
$ws.Range("Y32").Value = 200
$ws.Range("Z32").Value = 11.25
$ws.Range("AA32").Value = 11.75
$ws.Range("AB32").Value = 35
$ws.Range("AC32").Value = 250
$ws.Range("AE32").Value = 6.2
$ws.Range("AF32").Value = 5.2
$ws.Range("AG32").Value = 9.75
$ws.Range("AH32").Value = 6.4
$ws.Range("AI32").Value = 12
$ws.Range("AJ32").Value = 45
# Row 35
$ws.Range("K35").Value = 15
$ws.Range("AD35").Value = 301
$ws.Range("AE35").Value = 8.5
$ws.Range("AJ35").Value = 26
# Row 36
$ws.Range("G36").Value = 2.3
$ws.Range("H36").Value = 3.4
$ws.Range("I36").Value = 2.75
$ws.Range("P36").Value = 1.4
$ws.Range("Q36").Value = 2.75
$ws.Range("T36").Value = 8
$ws.Range("U36").Value = 11
$ws.Range("Y36").Value = 29
$ws.Range("AB36").Value = 15
$ws.Range("AC36").Value = 51
$ws.Range("AE36").Value = 9.5
# Row 37
$ws.Range("G37").Value = 2.63
$ws.Range("H37").Value = 3.25
$ws.Range("I37").Value = 2.7
$ws.Range("R37").Value = 1.67
$ws.Range("S37").Value = 2.1
$ws.Range("Y37").Value = 26
$ws.Range("AA37").Value = 6.5
$ws.Range("AB37").Value = 12
$ws.Range("AE37").Value = 10
# Row 39
$ws.Range("G39").Value = 2.5
$ws.Range("I39").Value = 2.9
$ws.Range("U39").Value = 11
$ws.Range("AG39").Value = 12
$ws.Range("AJ39").Value = 41
# Row 40
$ws.Range("R40").Value = 1.69
# Row 42
$ws.Range("K42").Value = 17
$ws.Range("R42").Value = 1.8
$ws.Range("S42").Value = 1.8
# Row 43
$ws.Range("R43").Value = 1.69
# Row 45
$ws.Range("G45").Value = 2.32
$ws.Range("I45").Value = 2.62
$ws.Range("R45").Value = 1.57
$ws.Range("S45").Value = 2.1
$ws.Range("T45").Value = 9.75
$ws.Range("W45").Value = 24
$ws.Range("Z45").Value = 13
$ws.Range("AH45").Value = 29
$ws.Range("AJ45").Value = 26
# Row 47
$ws.Range("G47").Value = 1.21
$ws.Range("I47").Value = 15
$ws.Range("L47").Value = 1.2
$ws.Range("M47").Value = 3.65
$ws.Range("N47").Value = 1.6
$ws.Range("O47").Value = 2.05
$ws.Range("R47").Value = 2.22
$ws.Range("S47").Value = 1.52
$ws.Range("T47").Value = 6.4
$ws.Range("U47").Value = 5.4
$ws.Range("W47").Value = 6.6
$ws.Range("X47").Value = 11.25
$ws.Range("Y47").Value = 35
$ws.Range("Z47").Value = 12
$ws.Range("AA47").Value = 11
$ws.Range("AB47").Value = 29
$ws.Range("AC47").Value = 150
$ws.Range("AE47").Value = 37
$ws.Range("AF47").Value = 150
$ws.Range("AG47").Value = 50
$ws.Range("AH47").Value = 800
$ws.Range("AI47").Value = 250
$ws.Range("AJ47").Value = 175
# Row 48
$ws.Range("G48").Value = 1.82
$ws.Range("H48").Value = 3.65
$ws.Range("I48").Value = 3.75
$ws.Range("L48").Value = 1.27
$ws.Range("M48").Value = 3.1
$ws.Range("N48").Value = 1.8
$ws.Range("O48").Value = 1.8
$ws.Range("S48").Value = 1.87
$ws.Range("T48").Value = 7.3
$ws.Range("U48").Value = 8.75
$ws.Range("W48").Value = 15
$ws.Range("X48").Value = 14.5
$ws.Range("Z48").Value = 11
$ws.Range("AA48").Value = 7
$ws.Range("AE48").Value = 11.25
$ws.Range("AF48").Value = 20
$ws.Range("AG48").Value = 13
$ws.Range("AH48").Value = 55
$ws.Range("AI48").Value = 35
$ws.Range("AJ48").Value = 40
# Row 49
$ws.Range("G49").Value = 1.25
$ws.Range("H49").Value = 4.65
$ws.Range("L49").Value = 1.22
$ws.Range("M49").Value = 3.45
$ws.Range("N49").Value = 1.65
$ws.Range("O49").Value = 1.98
$ws.Range("R49").Value = 2.15
$ws.Range("S49").Value = 1.55
$ws.Range("T49").Value = 6.2
$ws.Range("U49").Value = 5.6
$ws.Range("V49").Value = 8.75
$ws.Range("W49").Value = 7.3
$ws.Range("X49").Value = 11.25
$ws.Range("Y49").Value = 32
$ws.Range("AA49").Value = 9.75
$ws.Range("AB49").Value = 26
$ws.Range("AF49").Value = 110
$ws.Range("AI49").Value = 250
# Row 56
$ws.Range("G56").Value = 2.45
$ws.Range("I56").Value = 2.8
$ws.Range("L56").Value = 1.29
$ws.Range("M56").Value = 3.5
$ws.Range("R56").Value = 1.75
$ws.Range("S56").Value = 2
$ws.Range("T56").Value = 8.5
$ws.Range("U56").Value = 12
$ws.Range("W56").Value = 23
$ws.Range("X56").Value = 21
$ws.Range("Z56").Value = 10
$ws.Range("AD56").Value = 201
$ws.Range("AF56").Value = 13
$ws.Range("AH56").Value = 29
$ws.Range("AJ56").Value = 29
# Row 57
$ws.Range("N57").Value = 1.73
$ws.Range("O57").Value = 2.08
# Row 62
$ws.Range("G62").Value = 2.65
$ws.Range("H62").Value = 3.55
$ws.Range("I62").Value = 2.37
$ws.Range("P62").Value = 1.32
$ws.Range("Q62").Value = 3.1
$ws.Range("R62").Value = 1.53
$ws.Range("S62").Value = 2.32
$ws.Range("T62").Value = 11.5
$ws.Range("U62").Value = 15.5
$ws.Range("W62").Value = 30
$ws.Range("X62").Value = 19.5
$ws.Range("AA62").Value = 7.1
$ws.Range("AB62").Value = 11.75
$ws.Range("AC62").Value = 40
$ws.Range("AG62").Value = 9.25
$ws.Range("AH62").Value = 25
$ws.Range("AI62").Value = 17.5
